{"js": "// Replace \"Waktu Kampanye Rasi bintang Cygnus 2022:\" with\n// \"Waktu Kampanye 2022 untuk Rasi bintang Cygnus:\" everywhere it occurs\n// in the document body (4 occurrences), leaving the rest of the sentence\n// unchanged.\n\nconst searchText = \"Waktu Kampanye Rasi bintang Cygnus 2022:\";\nconst replaceText = \"Waktu Kampanye 2022 untuk Rasi bintang Cygnus:\";\n\nconst results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace \"Waktu Kampanye Rasi bintang Cygnus 2022:\" with\n# \"Waktu Kampanye 2022 untuk Rasi bintang Cygnus:\" everywhere it occurs\n# in the document body (4 occurrences), leaving the rest of the sentence\n# unchanged.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Waktu Kampanye Rasi bintang Cygnus 2022:\"\n$find.Replacement.Text = \"Waktu Kampanye 2022 untuk Rasi bintang Cygnus:\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n"}
